# [IMP] import_account_opening: Parter account from Excel file /Conto cliente/fornitore da file Excel
#
# Adds a "Ref" column between "Partita IVA" and "Dare"/"Avere" (Dare/Avere
# shift from F/G to G/H), and adds two extra example rows: a new customer
# row with a payment Ref of "RiBA" and a new supplier row, pushing the
# trailing "Banca" row down from row 4 to row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column F ("Ref"), shifting old Dare/Avere (F/G) to G/H ---
$ws.Columns("F").Insert()
$ws.Range("F1").Value = "Ref"

# --- Insert two new rows at position 4, shifting the old row 4 (Banca) to row 6 ---
$ws.Rows("4:5").Insert()

# --- Row 2: Prima Alpha S.p.A. (customer) ---
$ws.Range("B2").Value = "Prima Alpha S.p.A."
$ws.Range("E2").Value = "IT00115719999"

# --- Row 3: Notaio Decimo Jackson (supplier) ---
$ws.Range("B3").Value = "Notaio Decimo Jackson"
$ws.Range("E3").Value = "IT10242670015"

# --- New row 4: Latte Beta Due s.n.c. (customer, with a Ref of "RiBA") ---
$ws.Range("A4").Value = 152220
$ws.Range("B4").Value = "Latte Beta Due s.n.c."
$ws.Range("C4").Value = 1
$ws.Range("E4").Value = "IT02345670018"
$ws.Range("F4").Value = "RiBA"
$ws.Range("G4").Value = 150

# --- New row 5: Freie Universität Berlin (supplier) ---
$ws.Range("B5").Value = "Freie Universität Berlin"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "DE123456788"
$ws.Range("H5").Value = 200

$ws.Range("A6").Select()
